$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.448.52"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.564.77"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "288.83"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "0.3679"
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("D8").Value = "50.02"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "0.3373"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "1.139"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").Value = "0.07491"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "5.980"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "6.958"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "1.567.25"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "90.18"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "0.06745"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "6.365"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").Value = "16.19"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "12.03"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "22.458.25"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "2.391"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "2.623"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "19.76"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "149.28"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "5.053"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "124.23"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "1.742.35"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").Value = "  +5.62%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.019"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.164"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("D35").Value = "9.658"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "0.02461"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.341"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2272"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").Value = "0.06423"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "5.355"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "11.17"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "0.6172"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "3.763"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").Value = "0.5795"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "2.040"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "125.84"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "1.223"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "0.07322"
$ws.Range("E51").Value = "  +0.36%  "
